$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.791.08'
$ws.Range('E2').Value = '  +4.09%  '
$ws.Range('D3').Value = '3.479.24'
$ws.Range('E3').Value = '  +3.88%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '409.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +18.69%  '
$ws.Range('D7').Value = '3.473.80'
$ws.Range('E7').Value = '  +4.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.695'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.128'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +29.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.44%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.028.91'
$ws.Range('E13').Value = '  +3.85%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.142'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.54%  '
$ws.Range('D17').Value = '3.475.79'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '62.794.33'
$ws.Range('E18').Value = '  +4.51%  '
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000138'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +23.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '82.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '312.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '30.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.78'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.25%  '
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('B31').Value = 'LEO'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.120'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '43.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.17%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.41'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.126'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.53%  '
$ws.Range('E43').Value = '  +3.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '137.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.65'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.287'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('D50').Value = '2.209.22'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '3.822.33'
$ws.Range('E51').Value = '  +3.99%  '
